$d = $word.ActiveDocument

# Locate the title paragraph text that needs to be shortened and split into
# two runs with identical run formatting (bold title style).
$oldText = "APPLICATION FORM FOR OBTAINING ATTESTED FINGERPRINT SLIPS FOR FOREIGN ASSIGNMENTS"
$newText = "APPLICATION FORM FOR OBTAINING ATTESTED FINGERPRINT SLIPS"

$full = $d.Content.Text
$start = $full.IndexOf($oldText)
if ($start -lt 0) {
    throw "Could not find target title text in document."
}

# Replace the full title text with the shortened title (single run for now).
$titleRange = $d.Range($start, $start + $oldText.Length)
$titleRange.Text = $newText

# Split the new title into two runs - "...OBTAI" | "NING..." - matching the
# target edit, by forcing a run boundary at the split point. Toggling Bold
# off then back on creates a distinct run without leaving any stray
# formatting differences, since both halves end up with the same explicit
# Bold/BoldCs formatting as before.
$splitAt = 26
$secondRange = $d.Range($start + $splitAt, $start + $newText.Length)
$secondRange.Bold = 0
$secondRange.Bold = 1
